$d = $word.ActiveDocument
$p = $d.Paragraphs(4)
$r = $p.Range
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6D5C0788" w14:textId="77033172" w:rsidR="000C7DD9" w:rsidRDefault="000C7DD9" w:rsidP="00BA6DA6"><w:pPr><w:pStyle w:val="BasicParagraph"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="-72"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Optima-Regular"/></w:rPr></w:pPr><w:r><w:t>Está a participar numa campanha global para observar e registar as estrelas mais fracas visíveis como forma de medir a poluição luminosa num determinado local. Localizando e observando a  Constelação de Perseu no céu noturno e,  comparando-a com cartas estelares, pessoas de todo o mundo aprenderão  como as luzes da sua comunidade contribuem para a poluição luminosa. As suas contribuições para a base de dados on-line irão documentar a visibilidade do céu noturno em todo o mundo.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $r.InsertXML($xml)
